$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F29").Value = 76
$ws.Range("G29").Value = 3893.48
$ws.Range("F31").Value = 52
$ws.Range("G31").Value = 1385.28
$ws.Range("B34").Value = 61821.79
$ws.Range("F36").Value = 97
$ws.Range("G36").Value = 19086.69
$ws.Range("F41").Value = 225
$ws.Range("G41").Value = 43400.25
$ws.Range("F47").Value = 53
$ws.Range("G47").Value = 1340.37
$ws.Range("F52").Value = 60
$ws.Range("G52").Value = 3540
$ws.Range("F54").Value = 52
$ws.Range("G54").Value = 2319.2
$ws.Range("F61").Value = 242
$ws.Range("G61").Value = 63096.66
$ws.Range("F64").Value = 65
$ws.Range("G64").Value = 5168.8
$ws.Range("B66").Value = 219127.55
$ws.Range("B126").Value = 64196
$ws.Range("B127").Value = 65258
$ws.Range("F149").Value = 58
$ws.Range("G149").Value = 10514.82
$ws.Range("B155").Value = 39206.55
$ws.Range("F184").Value = 63
$ws.Range("G184").Value = 5166
$ws.Range("B193").Value = 69545.53
$ws.Range("F217").Value = 52
$ws.Range("G217").Value = 3863.6
$ws.Range("B218").Value = 84668.78
$ws.Range("F222").Value = 1068
$ws.Range("G222").Value = 19758
$ws.Range("B229").Value = 33627.37
$ws.Range("F276").Value = 14
$ws.Range("G276").Value = 910.84
$ws.Range("F284").Value = 174
$ws.Range("G284").Value = 8155.38
$ws.Range("B295").Value = 132155.9
$ws.Range("B304").Value = 55373
$ws.Range("E304").Value = 163.62
$ws.Range("F304").Value = -94
$ws.Range("G304").Value = -13562.32
$ws.Range("B305").Value = 63520
$ws.Range("E305").Value = 153.4
$ws.Range("F305").Value = 39
$ws.Range("G305").Value = 5626.92
$ws.Range("B306").Value = 57802
$ws.Range("E306").Value = 162.71
$ws.Range("F306").Value = -79
$ws.Range("G306").Value = -11334.92
$ws.Range("B307").Value = 63531
$ws.Range("E307").Value = 152.53
$ws.Range("F307").Value = 30
$ws.Range("G307").Value = 4304.4
$ws.Range("B308").Value = 55356
$ws.Range("E308").Value = 54.04
$ws.Range("F308").Value = -158
$ws.Range("G308").Value = -7527.12
$ws.Range("B309").Value = 63510
$ws.Range("E309").Value = 50.66
$ws.Range("F309").Value = 80
$ws.Range("G309").Value = 3811.2
$ws.Range("B317").Value = 63560
$ws.Range("E317").Value = 134.87
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 126.86
$ws.Range("B318").Value = 60325
$ws.Range("E318").Value = 151.57
$ws.Range("F318").Value = -102
$ws.Range("G318").Value = -12939.72
$ws.Range("F321").Value = 6
$ws.Range("G321").Value = 623.46
$ws.Range("F324").Value = 73
$ws.Range("G324").Value = 12507.09
$ws.Range("B328").Value = 1693.24
$ws.Range("F352").Value = 126
$ws.Range("G352").Value = 15738.66
$ws.Range("B356").Value = 80052.87
$ws.Range("F361").Value = 256
$ws.Range("G361").Value = 35991.04
$ws.Range("B363").Value = 81245.16
$ws.Range("F387").Value = 451
$ws.Range("G387").Value = 43566.6
$ws.Range("B389").Value = 60772.89
$ws.Range("F396").Value = 138
$ws.Range("G396").Value = 3516.24
$ws.Range("F397").Value = 79
$ws.Range("G397").Value = 2857.43
$ws.Range("F408").Value = 220
$ws.Range("G408").Value = 3487
$ws.Range("F410").Value = 75
$ws.Range("G410").Value = 18081
$ws.Range("F413").Value = 90
$ws.Range("G413").Value = 5181.3
$ws.Range("B417").Value = 178935.91
$ws.Range("F429").Value = 11
$ws.Range("G429").Value = 206.58
$ws.Range("F433").Value = 149
$ws.Range("G433").Value = 1436.36
$ws.Range("B438").Value = 27160.79
$ws.Range("F452").Value = 58
$ws.Range("G452").Value = 15730.76
$ws.Range("F454").Value = 80
$ws.Range("G454").Value = 22636.8
$ws.Range("F455").Value = 46
$ws.Range("G455").Value = 10222.58
$ws.Range("B458").Value = 102599.09
$ws.Range("B479").Value = 64810
$ws.Range("E479").Value = 291.22
$ws.Range("F479").Value = 0
$ws.Range("G479").Value = 0
$ws.Range("B480").Value = 53319
$ws.Range("E480").Value = 310.64
$ws.Range("F480").Value = -6
$ws.Range("G480").Value = -1643.52
$ws.Range("F511").Value = 261
$ws.Range("G511").Value = 26066.07
$ws.Range("F522").Value = 89
$ws.Range("G522").Value = 2370.96
$ws.Range("B525").Value = 132353.22
$ws.Range("F528").Value = 296
$ws.Range("G528").Value = 4694.56
$ws.Range("F529").Value = 128
$ws.Range("G529").Value = 4238.08
$ws.Range("F534").Value = 138
$ws.Range("G534").Value = 6038.88
$ws.Range("B535").Value = 26544.55
$ws.Range("F570").Value = 12
$ws.Range("G570").Value = 6418.56
$ws.Range("B573").Value = 30195.7
$ws.Range("F611").Value = 1
$ws.Range("G611").Value = 85.2
$ws.Range("F612").Value = 241
$ws.Range("G612").Value = 36248.81
$ws.Range("F617").Value = 35
$ws.Range("G617").Value = 1684.2
$ws.Range("F620").Value = 375
$ws.Range("G620").Value = 29471.25
$ws.Range("F622").Value = 493
$ws.Range("G622").Value = 50734.63
$ws.Range("F626").Value = 16
$ws.Range("G626").Value = 755.36
$ws.Range("B628").Value = 224393.98
$ws.Range("F674").Value = 952
$ws.Range("G674").Value = 155280.72
$ws.Range("B680").Value = 156293.27
$ws.Range("F706").Value = 124
$ws.Range("G706").Value = 4859.56
$ws.Range("F712").Value = 2
$ws.Range("G712").Value = 809.54
$ws.Range("B713").Value = 71793.24000000001
$ws.Range("B718").Value = 2973056.58
$ws.Range("B719").Value = 2973056.58
